$d = $word.ActiveDocument

$replacements = @(
    @('681÷7=97, 2', '326÷9=36, 2'),
    @('435÷2=217, 1', '355÷7=50, 5'),
    @('154÷4=38, 2', '304÷8=38, 0'),
    @('895÷6=149, 1', '366÷4=91, 2'),
    @('565÷8=70, 5', '567÷3=189, 0'),
    @('714÷2=357, 0', '809÷3=269, 2'),
    @('720÷9=80, 0', '468÷6=78, 0'),
    @('279÷7=39, 6', '855÷3=285, 0'),
    @('222÷4=55, 2', '955÷6=159, 1'),
    @('531÷7=75, 6', '880÷3=293, 1'),
    @('980÷7=140, 0', '203÷4=50, 3'),
    @('701÷2=350, 1', '685÷7=97, 6'),
    @('349÷6=58, 1', '865÷8=108, 1'),
    @('145÷7=20, 5', '693÷2=346, 1'),
    @('508÷4=127, 0', '315÷2=157, 1'),
    @('723÷7=103, 2', '660÷2=330, 0'),
    @('582÷5=116, 2', '764÷3=254, 2'),
    @('286÷8=35, 6', '798÷7=114, 0'),
    @('332÷9=36, 8', '880÷3=293, 1'),
    @('812÷9=90, 2', '114÷9=12, 6'),
    @('773÷9=85, 8', '367÷5=73, 2'),
    @('603÷8=75, 3', '495÷2=247, 1'),
    @('293÷3=97, 2', '938÷4=234, 2'),
    @('751÷9=83, 4', '342÷3=114, 0'),
    @('409÷7=58, 3', '405÷6=67, 3'),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()